$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'knee support sports',
    'knee support tape',
    'knee support tights',
    'knee support weightlifting',
    'knee support with hole',
    'knee support women',
    'knee support women workout',
    'knee support workout',
    'knee support workout women',
    'knee support x large',
    'knee support xl',
    'knee support xlarge',
    'knee support xs',
    'knee support yoga',
    'knee supports',
    'knee supports for basketball',
    'knee surgery games',
    'knee tape kinesiology',
    'knee tape support',
    'knee tights',
    'knee tights for women',
    'knee warmer cycling',
    'knee warmers for women',
    'knee warmers for women joint pain',
    'knee weightlifting support',
    'knee workout support',
    'knee wrap for running',
    'knee wraps for running',
    'knee wraps running',
    'knee yoga pants',
    'knees support women',
    'knees workout',
    'knew compression',
    'knit leggins',
    'kolorfull koala leggings',
    'kt ankle tape',
    'kt compression',
    'kt foot tape',
    'kt knee tape',
    'kt recovery ball',
    'kt recovery tape',
    'kt sport tape',
    'kyodan capri',
    'l a fitness',
    'l and m cigarette',
    'la apparel',
    'ladies athletic pants',
    'ladies basketball',
    'ladies boxing set',
    'ladies brown leggins',
    'ladies capri',
    'ladies capri leggings',
    'ladies capri pants',
    'ladies capries',
    'ladies capris',
    'ladies compression',
    'ladies compression capri leggings',
    'ladies compression hose',
    'ladies compression leggings',
    'ladies compression pants',
    'ladies compression shorts',
    'ladies compression shorts long',
    'ladies gym pants',
    'ladies gym tights',
    'ladies half pant',
    'ladies heavy tights',
    'ladies knit capris',
    'ladies leggings capri',
    'ladies leggings shorts',
    'ladies leggings tall',
    'ladies long leg panties',
    'ladies long length pants',
    'ladies long shorts',
    'ladies night game',
    'ladies riding tights',
    'ladies ski gear',
    'ladies sports leggings',
    'ladies stocking',
    'ladies stretch capri pants',
    'ladies summer capris',
    'ladies support hose',
    'ladies sweat pants petite',
    'ladies three quarter pants',
    'ladies tights',
    'ladies tights and leggings',
    'ladies tights leggings',
    'ladies underwear with support panel',
    'ladies white capri leggings',
    'ladies workout capri pants',
    'ladies workout tights',
    'lady apparel',
    'lady in the blue dress',
    'large hip pack',
    'large joints',
    'large tall running pants',
    'large tights',
    'large tights for women',
    'latisse leggings',
    'lay shirt exo',
    'leg air compression'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
